$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.03365660355876
$ws.Range("D2").Value = 1.040417923851197
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.049332497029906
$ws.Range("I2").Value = 1.035341824151251
$ws.Range("J2").Value = 1.03878013117511
$ws.Range("K2").Value = 1.043200329619352
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.052089849436717
$ws.Range("N2").Value = 1.016844914036938

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.034752827419709
$ws.Range("D3").Value = 1.041231920688309
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.050288928090319
$ws.Range("I3").Value = 1.035536020761288
$ws.Range("J3").Value = 1.039518416686534
$ws.Range("K3").Value = 1.043824736267942
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.052858132917898
$ws.Range("N3").Value = 1.01709382557943

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.035462174496043
$ws.Range("D4").Value = 1.041758513708321
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.050907988644583
$ws.Range("I4").Value = 1.035660306494792
$ws.Range("J4").Value = 1.039995609622565
$ws.Range("K4").Value = 1.044228002413924
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.053354828104865
$ws.Range("N4").Value = 1.017254599725252

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.035760388538861
$ws.Range("D5").Value = 1.04197986482503
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.051168285730087
$ws.Range("I5").Value = 1.035712227440821
$ws.Range("J5").Value = 1.04019609526365
$ws.Range("K5").Value = 1.044397351682575
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.053563534147192
$ws.Range("N5").Value = 1.017322120065882

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.035810460268548
$ws.Range("D6").Value = 1.042017028960716
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.051211993325459
$ws.Range("I6").Value = 1.035720925910506
$ws.Range("J6").Value = 1.04022975027854
$ws.Range("K6").Value = 1.044425775382913
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.053598570665775
$ws.Range("N6").Value = 1.017333452972635

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.035466159228265
$ws.Range("D7").Value = 1.041761471525731
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.050911466574517
$ws.Range("I7").Value = 1.035661001557138
$ws.Range("J7").Value = 1.039998289017749
$ws.Range("K7").Value = 1.044230265989762
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.053357617256698
$ws.Range("N7").Value = 1.017255502207028

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.03402707432553
$ws.Range("D8").Value = 1.040693041899858
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.049655688676398
$ws.Range("I8").Value = 1.035407737857265
$ws.Range("J8").Value = 1.039029747332365
$ws.Range("K8").Value = 1.04341150936168
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.052349584431731
$ws.Range("N8").Value = 1.016929094486636

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.031491330769042
$ws.Range("D9").Value = 1.0388094547351
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.047444288267368
$ws.Range("I9").Value = 1.034950952417049
$ws.Range("J9").Value = 1.037319018160429
$ws.Range("K9").Value = 1.041962894865205
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.05056997975612
$ws.Range("N9").Value = 1.016351718765317

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.0298008742181
$ws.Range("D10").Value = 1.037553166176707
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.04597101217686
$ws.Range("I10").Value = 1.034639378644349
$ws.Range("J10").Value = 1.036175819820159
$ws.Range("K10").Value = 1.040993224634334
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.049381359725967
$ws.Range("N10").Value = 1.015965324762559

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.029068888802404
$ws.Range("D11").Value = 1.0370090504436
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.045333306211328
$ws.Range("I11").Value = 1.034502793572001
$ws.Range("J11").Value = 1.035680158302826
$ws.Range("K11").Value = 1.040572417307385
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.048866151035911
$ws.Range("N11").Value = 1.015797662312871

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.028796995009158
$ws.Range("D12").Value = 1.036806921725546
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.045096468932327
$ws.Range("I12").Value = 1.034451808767472
$ws.Range("J12").Value = 1.035495949863131
$ws.Range("K12").Value = 1.040415970525158
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.048674700329185
$ws.Range("N12").Value = 1.015735332255992

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.028855317246384
$ws.Range("D13").Value = 1.036850279917832
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.045147269740479
$ws.Range("I13").Value = 1.034462756538898
$ws.Range("J13").Value = 1.03553546762176
$ws.Range("K13").Value = 1.040449535234369
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.048715770740392
$ws.Range("N13").Value = 1.015748704651866

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.029046414018527
$ws.Range("D14").Value = 1.036992342825113
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.045313728450518
$ws.Range("I14").Value = 1.03449858427241
$ws.Range("J14").Value = 1.035664933569389
$ws.Range("K14").Value = 1.040559488227259
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.048850327276939
$ws.Range("N14").Value = 1.015792511167055

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.029164154781081
$ws.Range("D15").Value = 1.037079869841832
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.045416293814308
$ws.Range("I15").Value = 1.034520625661174
$ws.Range("J15").Value = 1.035744688863793
$ws.Range("K15").Value = 1.040627215318869
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.048933221501004
$ws.Range("N15").Value = 1.015819494816524

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.029849453469988
$ws.Range("D16").Value = 1.037589274537779
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.046013339533867
$ws.Range("I16").Value = 1.034648408136846
$ws.Range("J16").Value = 1.036208701553042
$ws.Range("K16").Value = 1.041021132564893
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.049415541308217
$ws.Range("N16").Value = 1.015976444567515

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.030279320885986
$ws.Range("D17").Value = 1.037908774848606
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.046387912883178
$ws.Range("I17").Value = 1.034728115256506
$ws.Range("J17").Value = 1.03649959055838
$ws.Range("K17").Value = 1.041267976517769
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.049717946407395
$ws.Range("N17").Value = 1.016074800962513

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.030530054698378
$ws.Range("D18").Value = 1.038095120977551
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.046606417580619
$ws.Range("I18").Value = 1.034774445781127
$ws.Range("J18").Value = 1.036669198576869
$ws.Range("K18").Value = 1.041411866323969
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.04989428319701
$ws.Range("N18").Value = 1.01613213668732

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.030615548375992
$ws.Range("D19").Value = 1.038158657978274
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.046680925828802
$ws.Range("I19").Value = 1.034790215931917
$ws.Range("J19").Value = 1.036727019912168
$ws.Range("K19").Value = 1.041460913745376
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.049954400833027
$ws.Range("N19").Value = 1.016151680937742

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.030233200246141
$ws.Range("D20").Value = 1.03787449683104
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.046347722389151
$ws.Range("I20").Value = 1.034719580119043
$ws.Range("J20").Value = 1.036468387398879
$ws.Range("K20").Value = 1.04124150181606
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.049685506475168
$ws.Range("N20").Value = 1.016064251758036

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.028990140838951
$ws.Range("D21").Value = 1.036950509386921
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.045264709538396
$ws.Range("I21").Value = 1.034488040826588
$ws.Range("J21").Value = 1.035626811772515
$ws.Range("K21").Value = 1.040527113681089
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.048810705917637
$ws.Range("N21").Value = 1.015779612696398

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.028208569122662
$ws.Range("D22").Value = 1.036369446771681
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.044583979614691
$ws.Range("I22").Value = 1.034341010570956
$ws.Range("J22").Value = 1.035097114529385
$ws.Range("K22").Value = 1.0400771383627
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.048260225574171
$ws.Range("N22").Value = 1.015600343589526

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.028622896390575
$ws.Range("D23").Value = 1.036677489863368
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.044944828093025
$ws.Range("I23").Value = 1.034419091717902
$ws.Range("J23").Value = 1.035377970694028
$ws.Range("K23").Value = 1.040315755655013
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.048552089041011
$ws.Range("N23").Value = 1.015695406460849

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.030254040179531
$ws.Range("D24").Value = 1.037889985631781
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.046365882674108
$ws.Range("I24").Value = 1.034723437278591
$ws.Range("J24").Value = 1.036482486956857
$ws.Range("K24").Value = 1.04125346487245
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.04970016484162
$ws.Range("N24").Value = 1.016069018594149

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.032146871912103
$ws.Range("D25").Value = 1.039296508448631
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.048015815250673
$ws.Range("I25").Value = 1.035070285850322
$ws.Range("J25").Value = 1.037761760746619
$ws.Range("K25").Value = 1.042338089592698
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.051030442176041
$ws.Range("N25").Value = 1.016844914036938
